# Refresh the cryptocurrency price/volume snapshot (prices, % changes, and the
# Bittensor/PancakeSwap row order swap) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay a *text* cell (matches the sheet's existing
# inline-string cells) even when it looks numeric (e.g. "12.80", "0.999"), by using
# Excel's leading-apostrophe text-entry prefix, then restoring the default Normal
# style so no stray number-format/quote-prefix style is left on the cell.
function Set-TextValue([string]$addr, [string]$value) {
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = '65.669.22'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  +0.06%  '
Set-TextValue "D5" '600.46'
$ws.Range("E5").Value = '  -1.41%  '
Set-TextValue "D6" '156.79'
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +6.43%  '
Set-TextValue "D9" '0.130'
$ws.Range("E9").Value = '  +4.65%  '
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("E11").Value = '  -2.92%  '
$ws.Range("E12").Value = '  -0.21%  '
Set-TextValue "D13" '29.38'
$ws.Range("E13").Value = '  -2.86%  '
$ws.Range("E14").Value = '  -1.39%  '
$ws.Range("D15").Value = '3.157.56'
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").Value = '65.521.30'
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").Value = '2.682.59'
$ws.Range("E17").Value = '  -0.50%  '
Set-TextValue "D18" '12.80'
$ws.Range("E18").Value = '  +1.06%  '
$ws.Range("E19").Value = '  -1.88%  '
$ws.Range("E20").Value = '  -0.05%  '
Set-TextValue "D21" '351.62'
$ws.Range("E21").Value = '  -2.30%  '
$ws.Range("E22").Value = '  -0.01%  '
Set-TextValue "D23" '69.63'
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("E24").Value = '  +4.97%  '
Set-TextValue "D25" '9.66'
$ws.Range("E25").Value = '  -1.49%  '
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("E27").Value = '  -0.79%  '
$ws.Range("E28").Value = '  -5.45%  '
Set-TextValue "D29" '8.10'
$ws.Range("E29").Value = '  -1.69%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D31" '2.15'
$ws.Range("E31").Value = '  -2.38%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D32" '531.49'
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("E33").Value = '  -2.05%  '
$ws.Range("E34").Value = '  -2.51%  '
Set-TextValue "D35" '5.49'
$ws.Range("E35").Value = '  +0.84%  '
Set-TextValue "D36" '0.425'
$ws.Range("E37").Value = '  -1.37%  '
Set-TextValue "D39" '158.06'
$ws.Range("E39").Value = '  -3.08%  '
$ws.Range("E40").Value = '  -2.41%  '
Set-TextValue "D41" '0.999'
$ws.Range("E41").Value = '  +0.02%  '
Set-TextValue "D42" '164.40'
$ws.Range("E42").Value = '  -2.60%  '
$ws.Range("E43").Value = '  -0.42%  '
Set-TextValue "D44" '2.32'
$ws.Range("E44").Value = '  +2.11%  '
$ws.Range("E45").Value = '  -0.32%  '
Set-TextValue "D46" '22.84'
$ws.Range("E46").Value = '  -2.77%  '
$ws.Range("E47").Value = '  -2.07%  '
$ws.Range("E48").Value = '  -2.35%  '
$ws.Range("D49").Value = '0.0₆0261'
$ws.Range("E49").Value = '  +15.61%  '
$ws.Range("E50").Value = '  +2.89%  '
Set-TextValue "D51" '20.10'
$ws.Range("E51").Value = '  -5.31%  '
